# "Saving models and results feature"
#
# Slide 1 title placeholder: the old "Chuong 03 / GIOI THIEU / FEATURE
# ENGINEERING" (multi-run, multi-line) title is retyped as a single line,
# "HANDWRITTEN DIGIT IDENTIFICATION" - exactly as if the author selected
# all but the trailing character, deleted it (collapsing every run/line
# break down to the formatting of the last run, which was not bold) and
# then typed the replacement text, leaving the paragraph's endParaRPr
# untouched.
$p = $ppt.ActivePresentation

$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Title
$range1 = $title1.TextFrame.TextRange
$lastChar = $range1.Length
$range1.Characters(1, $lastChar - 1).Delete()
$range1.Text = "HANDWRITTEN DIGIT IDENTIFICATION"

# Slide 2 title placeholder: text content is unchanged ("Handwritten
# digit identification") - just touched/re-saved.
$slide2 = $p.Slides.Item(2)
$title2 = $slide2.Shapes.Title
$title2.TextFrame.TextRange.Text = "Handwritten digit identification"
